$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 5 (pushes the old "Source" row, currently
#    row 5, down to row 6) so we end up with two data rows (4 and 5) plus
#    the source row at 6.
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).Insert()

# ---------------------------------------------------------------------------
# 2. Row 1 - title, now spans A1:I1
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Merge()
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Vani Municipality"
$ws.Range("A1:I1").Font.Name = "Arial"
$ws.Range("A1:I1").Font.Size = 11
$ws.Range("A1:I1").Font.Bold = $true
$ws.Range("A1:I1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1:I1").VerticalAlignment = -4108     # xlCenter
$ws.Range("A1:I1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 51

# ---------------------------------------------------------------------------
# 3. Row 2 - "(End of year, persons)" - unchanged content, keep formatting
# ---------------------------------------------------------------------------
# (left as-is)

# ---------------------------------------------------------------------------
# 4. Row 3 - A3 font becomes Sylfaen
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11

# ---------------------------------------------------------------------------
# 5. Row 4 - "family with disabilities Persons " + first data series
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none
$ws.Rows.Item(4).RowHeight = 24.75

$row4vals = 877,849,796,796,767,758,735,725
$cols = "B","C","D","E","F","G","H","I"
for ($i = 0; $i -lt 8; $i++) {
    $cell = $ws.Range($cols[$i] + "4")
    $cell.Value = $row4vals[$i]
    $cell.Borders.Item(8).LineStyle = -4142   # xlEdgeTop -> none
    $cell.Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none
    $cell.HorizontalAlignment = -4131          # xlGeneral
}

# ---------------------------------------------------------------------------
# 6. Row 5 (new row) - "disabilities Persons " + second data series
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10
$ws.Range("A5").HorizontalAlignment = -4131   # xlLeft ... set precisely below
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true
$ws.Range("A5").Borders.Item(8).LineStyle = -4142  # top none
$ws.Range("A5").Borders.Item(9).LineStyle = 1      # bottom thin
$ws.Rows.Item(5).RowHeight = 21

$row5vals = 997,961,898,900,866,860,833,827
for ($i = 0; $i -lt 8; $i++) {
    $cell = $ws.Range($cols[$i] + "5")
    $cell.Value = $row5vals[$i]
    $cell.NumberFormat = "#\ ##0"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.HorizontalAlignment = -4131
    if ($cols[$i] -eq "I") {
        $cell.Borders.Item(9).LineStyle = 1    # bottom thin for I5
    } else {
        $cell.Borders.Item(9).LineStyle = -4142
    }
    $cell.Borders.Item(8).LineStyle = -4142
}

# ---------------------------------------------------------------------------
# 7. Row 6 - Source row (shifted down); remove the top border on A6
# ---------------------------------------------------------------------------
$ws.Range("A6").Borders.Item(8).LineStyle = -4142   # xlEdgeTop -> none
$ws.Rows.Item(6).RowHeight = 27.75

# ---------------------------------------------------------------------------
# 8. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20

$wb.Save()
